$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Fix the "stuff after this line >>>" paragraph: drop the trailing ">"
#    so it reads "...this line >>". We rewrite the paragraph's three runs
#    verbatim (via WordOpenXML) so only the third run's text changes and
#    the existing run/paragraph identity (rsids, paraId) is preserved -
#    a plain Find/Replace or Range.Text edit would otherwise cause the
#    engine to merge the unformatted runs together.
# ---------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith(">>>")) {
        $target = $p
        break
    }
}

$pkgXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="647D9B56" w14:textId="739F44DD" w:rsidR="00091943" w:rsidRDefault="00091943">
<w:r><w:t>&gt;&gt;</w:t></w:r>
<w:r w:rsidR="00432DF3"><w:t>&gt; your</w:t></w:r>
<w:r><w:t xml:space="preserve"> stuff after this line &gt;&gt;</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$target.Range.InsertXML($pkgXml)

# ---------------------------------------------------------------------
# 2. Insert two new paragraphs right after that paragraph:
#      a) a red (FF0000) colored paragraph of commentary text
#      b) a following blank paragraph
# ---------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith(">>")) {
        $target = $p
        break
    }
}

$insertionPoint = $target.Range
$insertionPoint.Collapse(0)             # wdCollapseEnd
$insertionPoint.InsertParagraphAfter()  # creates the (still empty) red-text paragraph
$redPara = $target.Next()
$redPara.Range.InsertParagraphAfter()   # creates the trailing blank paragraph

# Fill the red-text paragraph with its text + red run/paragraph-mark formatting.
$redPara = $target.Next()
$redXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr>
<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Git hub is so cool. I love being able to manage different version of my software when collaborating with other people regarding personal projects and work.  Hope everyone is doing well with covid and all that Jazz, stay safe people.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$redPara.Range.InsertXML($redXml)

# Tidy the trailing blank paragraph so it serializes as a bare <w:p/>
# rather than <w:p><w:r/></w:p>.
$blankPara = $redPara.Next()
$blankRng = $blankPara.Range
$blankRng.Collapse(1)  # wdCollapseStart
$blankXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$blankRng.InsertXML($blankXml)
